# Revision of "Loading your Data" documentation sample workbook.
# - Header labels get a type-prefix added (mD#, mS#, c#, i#) so the
#   docs tool can recognize column roles.
# - The window position / last-selected-cell bookkeeping is refreshed
#   to reflect where the author left off while editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (add semantic-type prefixes) -------------
$ws.Range("A1").Value = "mD#function"
$ws.Range("B1").Value = "mS#gene"
$ws.Range("E1").Value = "c#heat 0"
$ws.Range("F1").Value = "i#heat 10"
$ws.Range("G1").Value = "i#heat 20"

# --- Move the active selection on the sheet ---------------------------
$ws.Range("G2").Select()

# --- Reposition the workbook window on screen --------------------------
$win = $excel.ActiveWindow
$win.Left = 7240
$win.Top = 11740
